$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '20.499.78'
Set-TextValue $ws.Range("E2") '  +2.91%  '
Set-TextValue $ws.Range("D3") '1.468.76'
Set-TextValue $ws.Range("E3") '  +4.17%  '
Set-TextValue $ws.Range("E4") '  +0.71%  '
Set-TextValue $ws.Range("D5") '280.92'
Set-TextValue $ws.Range("E5") '  +3.06%  '
Set-TextValue $ws.Range("D6") '0.9006'
Set-TextValue $ws.Range("E6") '  -10.22%  '
Set-TextValue $ws.Range("D7") '0.3721'
Set-TextValue $ws.Range("E7") '  +0.57%  '
Set-TextValue $ws.Range("D8") '0.3197'
Set-TextValue $ws.Range("E8") '  +4.44%  '
Set-TextValue $ws.Range("D9") '40.32'
Set-TextValue $ws.Range("E9") '  +3.25%  '
Set-TextValue $ws.Range("D10") '1.055'
Set-TextValue $ws.Range("E10") '  +6.26%  '
Set-TextValue $ws.Range("D11") '0.06678'
Set-TextValue $ws.Range("E11") '  +2.03%  '
Set-TextValue $ws.Range("D12") '1.004'
Set-TextValue $ws.Range("E12") '  +0.04%  '
Set-TextValue $ws.Range("D13") '5.566'
Set-TextValue $ws.Range("E13") '  +2.90%  '
Set-TextValue $ws.Range("D14") '18.08'
Set-TextValue $ws.Range("E14") '  +7.75%  '
Set-TextValue $ws.Range("D15") '6.236'
Set-TextValue $ws.Range("E15") '  +1.41%  '
Set-TextValue $ws.Range("D16") '1.474.33'
Set-TextValue $ws.Range("E16") '  +4.14%  '
Set-TextValue $ws.Range("D17") '0.00001037'
Set-TextValue $ws.Range("E17") '  +3.31%  '
Set-TextValue $ws.Range("D18") '0.05708'
Set-TextValue $ws.Range("E18") '  -0.43%  '
Set-TextValue $ws.Range("D19") '0.9002'
Set-TextValue $ws.Range("E19") '  -10.22%  '
Set-TextValue $ws.Range("D20") '71.28'
Set-TextValue $ws.Range("E20") '  -3.29%  '
Set-TextValue $ws.Range("D21") '5.709'
Set-TextValue $ws.Range("E21") '  +2.59%  '
Set-TextValue $ws.Range("D22") '14.73'
Set-TextValue $ws.Range("E22") '  +2.57%  '
Set-TextValue $ws.Range("D23") '11.26'
Set-TextValue $ws.Range("E23") '  +4.44%  '
Set-TextValue $ws.Range("D24") '2.298'
Set-TextValue $ws.Range("E24") '  -1.15%  '
Set-TextValue $ws.Range("D25") '20.798.79'
Set-TextValue $ws.Range("E25") '  +4.34%  '
Set-TextValue $ws.Range("D26") '2.296'
Set-TextValue $ws.Range("E26") '  +2.03%  '
Set-TextValue $ws.Range("D27") '138.14'
Set-TextValue $ws.Range("E27") '  -0.54%  '
Set-TextValue $ws.Range("E28") '  +3.77%  '
Set-TextValue $ws.Range("D29") '1.639.09'
Set-TextValue $ws.Range("E29") '  +4.11%  '
Set-TextValue $ws.Range("D30") '113.70'
Set-TextValue $ws.Range("E30") '  +4.29%  '
Set-TextValue $ws.Range("D31") '3.958'
Set-TextValue $ws.Range("E31") '  +2.73%  '
Set-TextValue $ws.Range("D32") '5.250'
Set-TextValue $ws.Range("E32") '  -1.69%  '
Set-TextValue $ws.Range("D33") '0.8521'
Set-TextValue $ws.Range("E33") '  +0.04%  '
Set-TextValue $ws.Range("D34") '0.07829'
Set-TextValue $ws.Range("E34") '  +1.89%  '
Set-TextValue $ws.Range("B35") 'Hedera'
Set-TextValue $ws.Range("C35") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D35") '0.06165'
Set-TextValue $ws.Range("E35") '  +7.57%  '
Set-TextValue $ws.Range("B36") 'WEMIXTOKEN'
Set-TextValue $ws.Range("C36") 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D36") '1.513'
Set-TextValue $ws.Range("E36") '  +18.17%  '
Set-TextValue $ws.Range("E37") '  +2.62%  '
Set-TextValue $ws.Range("D38") '1.175'
Set-TextValue $ws.Range("E38") '  +11.52%  '
Set-TextValue $ws.Range("D39") '10.66'
Set-TextValue $ws.Range("E39") '  +2.98%  '
Set-TextValue $ws.Range("D40") '0.02062'
Set-TextValue $ws.Range("E40") '  +1.79%  '
Set-TextValue $ws.Range("D41") '0.1889'
Set-TextValue $ws.Range("E41") '  -1.07%  '
Set-TextValue $ws.Range("D42") '0.9260'
Set-TextValue $ws.Range("E42") '  -7.61%  '
Set-TextValue $ws.Range("D43") '7.307'
Set-TextValue $ws.Range("E43") '  -12.83%  '
Set-TextValue $ws.Range("D44") '0.5394'
Set-TextValue $ws.Range("E44") '  +2.41%  '
Set-TextValue $ws.Range("D45") '3.590'
Set-TextValue $ws.Range("E45") '  +1.90%  '
Set-TextValue $ws.Range("D46") '12.48'
Set-TextValue $ws.Range("E46") '  +1.55%  '
Set-TextValue $ws.Range("D47") '124.01'
Set-TextValue $ws.Range("E47") '  +13.70%  '
Set-TextValue $ws.Range("D48") '0.5305'
Set-TextValue $ws.Range("E48") '  +4.28%  '
Set-TextValue $ws.Range("D49") '1.832'
Set-TextValue $ws.Range("E49") '  +2.40%  '
Set-TextValue $ws.Range("D50") '0.06450'
Set-TextValue $ws.Range("E50") '  +5.27%  '
Set-TextValue $ws.Range("E51") '  +0.01%  '
